# feat: add 2022-Q1 data
#
# 1. A new worksheet "2022-Q1" is inserted right before the "总计" (Total)
#    summary worksheet, containing the per-fund holding detail for the new
#    quarter.
# 2. The "总计" worksheet gets a new first data row summarizing the
#    2022-Q1 quarter (holding count = 3, holding value = 0.56 billion).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: apply the "header / index-column" look used throughout this
# workbook - bold font, thin border, centered & top-aligned text.
# ---------------------------------------------------------------------
function Format-HeaderCell($rng) {
    $rng.Font.Bold = $true
    $rng.HorizontalAlignment = -4108   # xlCenter
    $rng.VerticalAlignment = -4160     # xlTop
    $rng.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# Locate the existing "总计" sheet and remember its current data so we
# can rebuild it (with the extra 2022-Q1 row) after re-creating it in
# its new position.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

$oldRows = @(
    @("2021-Q4", 5, 0.9),
    @("2021-Q3", 3, 0.61),
    @("2021-Q2", 3, 0.54),
    @("2021-Q1", 4, 0.51),
    @("2020-Q4", 3, 0.4)
)

# Remove the old "总计" sheet - it will be re-added (after the new
# "2022-Q1" sheet) further down so that the internal sheet ordering /
# identifiers match a freshly appended sheet.
[void]$totalSheet.Delete()

# ---------------------------------------------------------------------
# Add the new "2022-Q1" worksheet right after "2021-Q4" (i.e. at the
# end of the workbook, since "总计" was just removed).
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$q1.Name = "2022-Q1"

$q1headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $q1headers.Length; $c++) {
    $cell = $q1.Cells.Item(1, 2 + $c)
    Format-HeaderCell $cell
    $cell.Value = $q1headers[$c]
}

$q1data = @(
    @("002423", "华宝兴业标普美国消费(QDII-LOF)美元", "3.62", "94.37", "7.14", "0.2585", 3),
    @("162415", "华宝标普美国消费(QDII-LOF)人民币A",   "3.62", "94.37", "7.14", "0.2585", 3),
    @("009975", "华宝标普美国消费(QDII-LOF)人民币C",   "0.61", "94.37", "7.14", "0.0436", 3)
)

for ($i = 0; $i -lt $q1data.Length; $i++) {
    $row = 2 + $i
    $rec = $q1data[$i]

    $idxCell = $q1.Cells.Item($row, 1)
    Format-HeaderCell $idxCell
    $idxCell.Value = $i

    $q1.Cells.Item($row, 2).Value = "'" + $rec[0]
    $q1.Cells.Item($row, 3).Value = $rec[1]
    $q1.Cells.Item($row, 4).Value = "'" + $rec[2]
    $q1.Cells.Item($row, 5).Value = "'" + $rec[3]
    $q1.Cells.Item($row, 6).Value = "'" + $rec[4]
    $q1.Cells.Item($row, 7).Value = "'" + $rec[5]
    $q1.Cells.Item($row, 8).Value = $rec[6]
}

[void]$q1.Range("A1").Select()

# ---------------------------------------------------------------------
# Re-add the "总计" worksheet at the end, with the original rows plus
# the new 2022-Q1 summary row on top.
# ---------------------------------------------------------------------
$lastIndex = $wb.Worksheets.Count
$total = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($lastIndex))
$total.Name = "总计"

$totalHeaders = @("日期", "持有数量(只)", "持有市值(亿元)")
for ($c = 0; $c -lt $totalHeaders.Length; $c++) {
    $cell = $total.Cells.Item(1, 2 + $c)
    Format-HeaderCell $cell
    $cell.Value = $totalHeaders[$c]
}

$allRows = New-Object System.Collections.ArrayList
[void]$allRows.Add(@("2022-Q1", 3, 0.56))
foreach ($r in $oldRows) {
    [void]$allRows.Add($r)
}

for ($i = 0; $i -lt $allRows.Count; $i++) {
    $row = 2 + $i
    $rec = $allRows[$i]

    $idxCell = $total.Cells.Item($row, 1)
    Format-HeaderCell $idxCell
    $idxCell.Value = $i

    $total.Cells.Item($row, 2).Value = $rec[0]
    $total.Cells.Item($row, 3).Value = $rec[1]
    $total.Cells.Item($row, 4).Value = $rec[2]
}

[void]$total.Range("A1").Select()
